$d = $word.ActiveDocument

# Paragraph 30 (1-indexed): "...and a customer. – NOT DONE"
# Remove the word "NOT " so the run splits into "–" / " " / "DONE"
$p1 = $d.Paragraphs.Item(30)
$r1 = $p1.Range
$r1.Find.Execute("NOT ", $false, $false, $false, $false, $false, $true, 1, $false, "", 2)

# Paragraph 40 (1-indexed): "...by the website. – NOT DONE"
# Remove " NOT" (leading space + NOT) so the run text becomes "–DONE"
$p2 = $d.Paragraphs.Item(40)
$r2 = $p2.Range
$r2.Find.Execute(" NOT DONE", $false, $false, $false, $false, $false, $true, 1, $false, "DONE", 2)
